# Generate Report for Handback
# - Marks zh-cn / de-de as handed back (status text + handback datetime)
# - Populates the "Latest Target File" / "Latest Handback File" columns
#   for both language sheets, with a hyperlink on the target-file cell
# - Widens a few report columns that now hold longer text

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetFile = "e7a4e73a-68d0-48e9-a82c-b6de2683f2bf.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/967a4cbe24647a0d87880306d811f44c3ef5a303/e2e/e7a4e73a-68d0-48e9-a82c-b6de2683f2bf.md"
$hyperlinkColor = 15570276  # OLE BGR for #6495ED (cornflower blue), matches the workbook's existing HyperLink style

# ---------------------------------------------------------------------
# Overview sheet: both language status cells flip to "handed back"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# Overview columns E (zh-cn) and F (de-de) grow to fit the longer status text
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText

$zhcn.Range("I2").Value = $targetFile
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetUrl, "", "", $targetFile)
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = $hyperlinkColor

$zhcn.Range("J2").Value = "e7a4e73a-68d0-48e9-a82c-b6de2683f2bf.324a880f8954ffbdc76951ba4765a752eb6a73a0.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-28 08:57:32"

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText

$dede.Range("I2").Value = $targetFile
$dede.Hyperlinks.Add($dede.Range("I2"), $targetUrl, "", "", $targetFile)
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = $hyperlinkColor

$dede.Range("J2").Value = "e7a4e73a-68d0-48e9-a82c-b6de2683f2bf.324a880f8954ffbdc76951ba4765a752eb6a73a0.de-de.xlf"
$dede.Range("K2").Value = "2016-08-28 08:57:39"

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Handback report generated."
